$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.105.72"
$ws.Range("E2").Value = "  +3.72%  "
$ws.Range("D3").Value = "2.420.34"
$ws.Range("E3").Value = "  +3.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "553.37"
$ws.Range("E5").Value = "  +2.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.82"
$ws.Range("E6").Value = "  +2.80%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.577"
$ws.Range("E8").Value = "  +2.66%  "
$ws.Range("E9").Value = "  +3.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.77"
$ws.Range("E10").Value = "  +4.91%  "
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "24.88"
$ws.Range("E13").Value = "  +4.50%  "
$ws.Range("D14").Value = "2.849.99"
$ws.Range("E14").Value = "  +3.18%  "
$ws.Range("D15").Value = "59.999.07"
$ws.Range("E16").Value = "  +3.19%  "
$ws.Range("D17").Value = "2.408.34"
$ws.Range("E17").Value = "  +2.80%  "
$ws.Range("E18").Value = "  +6.29%  "
$ws.Range("E19").Value = "  +2.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "333.11"
$ws.Range("E20").Value = "  +1.36%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.32"
$ws.Range("E23").Value = "  +3.78%  "
$ws.Range("E24").Value = "  +3.81%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.60"
$ws.Range("E25").Value = "  +3.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("D28").Value = "0.0₃0784"
$ws.Range("E28").Value = "  +6.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "169.38"
$ws.Range("E31").Value = "  -0.56%  "
$ws.Range("E32").Value = "  +2.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.68"
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.22"
$ws.Range("E37").Value = "  +1.21%  "
$ws.Range("E38").Value = "  +0.59%  "
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "39.45"
$ws.Range("E39").Value = "  +0.95%  "
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.418"
$ws.Range("E40").Value = "  +10.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "313.34"
$ws.Range("E41").Value = "  +8.57%  "
$ws.Range("E42").Value = "  +1.68%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "139.20"
$ws.Range("E43").Value = "  -1.46%  "
$ws.Range("E44").Value = "  +1.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0520"
$ws.Range("E45").Value = "  +2.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.53"
$ws.Range("E46").Value = "  +2.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.413"
$ws.Range("E47").Value = "  +8.45%  "
$ws.Range("E48").Value = "  +1.46%  "
$ws.Range("E49").Value = "  +1.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "17.72"
$ws.Range("E50").Value = "  +2.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.06"
$ws.Range("E51").Value = "  -0.19%  "
